$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 14 with date, depth, and a new note
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A14").Value = (Get-Date -Year 2023 -Month 5 -Day 12 -Hour 10 -Minute 24 -Second 0)
$ws.Range("B14").Value = 2.24
$ws.Range("C14").Value = "Brush may have been removed from grate"

# Fix the typo in the existing "Rained this morning" note (C12)
$ws.Range("C12").Value = "Rained this morning (~0.23"" according to CoCoRaHS)"
